# edit.ps1
# Applies the BasicItem_Common.xlsx data-fix edit:
#  - Clears the stray "{}" placeholder values in column H (ConsumeItem) for the
#    data rows (rows 4-30) on Sheet1, which in turn lets Excel drop the now-unused
#    duplicate "{}" shared-string table entries when the workbook is saved.
#  - Fixes two bad data values on row 7: the item Id (column A) and quality
#    (column D) were mistyped as 800011/12 and should be 5/4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear out column H ("ConsumeItem") placeholder values for the data rows.
$ws.Range("H4:H30").ClearContents()

# Correct the Id and Quality values on row 7.
$ws.Range("A7").Value = 5
$ws.Range("D7").Value = 4
